$wb = $excel.ActiveWorkbook

# Data QA fix: column A holds a Y/N flag per report row. Every "N" should
# actually read "Y" (the diff drops the now-unused "N" shared string and
# repoints every row at the existing "Y" string).
$lastRows = @{1 = 121; 2 = 41; 3 = 71; 4 = 81; 5 = 21; 6 = 64}

foreach ($sheetIndex in 1..6) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $lastRow = $lastRows[$sheetIndex]
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        if ($cell.Text -eq "N") {
            $cell.Value = "Y"
        }
    }
}

# Restore each sheet's scroll position / selection the way the author last
# left it, and move the active tab to "Class Status".
$ws1 = $wb.Worksheets.Item(1)
$ws1.Application.ActiveWindow.ScrollRow = 110
$ws1.Range("D117").Select()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Application.ActiveWindow.ScrollRow = 32
$ws2.Range("A2:A41").Select()

$ws3 = $wb.Worksheets.Item(3)
$ws3.Application.ActiveWindow.ScrollRow = 59
$ws3.Range("A2:A71").Select()

$ws4 = $wb.Worksheets.Item(4)
$ws4.Application.ActiveWindow.ScrollRow = 72
$ws4.Range("A2:A81").Select()

$ws5 = $wb.Worksheets.Item(5)
$ws5.Application.ActiveWindow.ScrollRow = 18
$ws5.Range("E21").Select()

$ws6 = $wb.Worksheets.Item(6)
$ws6.Application.ActiveWindow.ScrollRow = 59
$ws6.Range("E61").Select()

$ws6.Activate()
